# Update the Team Expense Split workbook:
#  - "Shared Expenses" sheet: update Amount Spent values for Bergi/Nils/Adrien/Jason,
#    with Jason's amount now computed by a formula (7.04+21.83), and move the
#    selected cell.
#  - "Direct Expenses" sheet: move the selected cell.
# Dependent formulas on "Balances" and "Owes Matrix" recalc automatically.

$wb = $excel.ActiveWorkbook

$shared = $wb.Worksheets.Item("Shared Expenses")
$shared.Range("B2").Value = 29.88
$shared.Range("B3").Value = 0
$shared.Range("B4").Value = 0
$shared.Range("B5").Formula = "=7.04+21.83"

$shared.Range("E12").Select()

$direct = $wb.Worksheets.Item("Direct Expenses")
$direct.Range("G24").Select()

$shared.Activate()

$excel.CalculateFull()
